# B2C hotel booking flow and modification on config file for all env
#
# 1) Extend the "FlightEmployeeData" sheet with new flight / user columns
# 2) Add a brand-new "B2CHotelBookingData" sheet with hotel-booking data

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. FlightEmployeeData additions
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("FlightEmployeeData")

# Style "donor" cells already present on the workbook so that the new
# cells reuse the existing style entries instead of creating duplicates.
$styleHeader = $ws4.Range("A1")   # plain header style
$styleText   = $ws4.Range("D2")   # text/number-as-text style

# Widen column H (the new 8th column) the same as column G (27 chars wide)
$ws4.Columns.Item(8).ColumnWidth = 26.1666666667

# --- Row 1 (headers) ---------------------------------------------------
$cells1 = @("H1","I1","J1","K1","L1","M1","N1","O1")
$values1 = @("Flight Name","Price Type","AddOns","City","Username","usernumber","userdate","usermob")
for ($i = 0; $i -lt $cells1.Length; $i++) {
    $c = $ws4.Range($cells1[$i])
    $styleHeader.Copy()
    $c.PasteSpecial(-4122)
    $c.Value = $values1[$i]
}

# --- Row 2 (data), processed left-to-right so new shared strings are
#     appended in the same relative order the reference workbook used.
$wrapDonor = $wb.Worksheets.Item("Sheet1").Range("A3")

# H2 -> text style (same as D2)
$h2 = $ws4.Range("H2")
$styleText.Copy()
$h2.PasteSpecial(-4122)
$h2.Value = "SpiceJet"

# I2 -> brand-new style: same font/fill as the wrap-text style used on
# Sheet1!A3 (fontId 2 / fillId 2) but without wrap-text.
$i2 = $ws4.Range("I2")
$wrapDonor.Copy()
$i2.PasteSpecial(-4122)
$i2.WrapText = $false
$i2.Value = "Regular"

# J2 -> plain style
$j2 = $ws4.Range("J2")
$styleHeader.Copy()
$j2.PasteSpecial(-4122)
$j2.Value = "No"

# K2 -> plain style (value already an existing shared string: "Delhi")
$k2 = $ws4.Range("K2")
$styleHeader.Copy()
$k2.PasteSpecial(-4122)
$k2.Value = "Delhi"

# L2 -> plain style
$l2 = $ws4.Range("L2")
$styleHeader.Copy()
$l2.PasteSpecial(-4122)
$l2.Value = "Vikas1234"

# M2 -> text style
$m2 = $ws4.Range("M2")
$styleText.Copy()
$m2.PasteSpecial(-4122)
$m2.Value = "vikas1234"

# N2 -> text style
$n2 = $ws4.Range("N2")
$styleText.Copy()
$n2.PasteSpecial(-4122)
$n2.Value = "Jan 2024 20"

# O2 -> text style (value already an existing shared string: "9319580173")
$o2 = $ws4.Range("O2")
$styleText.Copy()
$o2.PasteSpecial(-4122)
$o2.Value = "9319580173"

# ---------------------------------------------------------------------
# 2. New sheet: B2CHotelBookingData
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = "B2CHotelBookingData"

$ws5.Columns.Item(1).ColumnWidth = 18.8333333333
$ws5.Columns.Item(2).ColumnWidth = 16.1666666667
$ws5.Columns.Item(3).ColumnWidth = 15.0

# header row uses the plain style (same as FlightEmployeeData A1)
$cells1b = @("A1","B1","C1","D1")
$values1b = @("City","CheckInDate","CheckOutdate","Guest")
for ($i = 0; $i -lt $cells1b.Length; $i++) {
    $c = $ws5.Range($cells1b[$i])
    $styleHeader.Copy()
    $c.PasteSpecial(-4122)
    $c.Value = $values1b[$i]
}

# A2 -> plain style
$a2 = $ws5.Range("A2")
$styleHeader.Copy()
$a2.PasteSpecial(-4122)
$a2.Value = "New Delhi"

# B2:D2 -> text style (same as D2 on FlightEmployeeData)
$cells2b = @("B2","C2","D2")
$values2b = @("Jan 2024 22","Jan 2024 25","1 Guest")
for ($i = 0; $i -lt $cells2b.Length; $i++) {
    $c = $ws5.Range($cells2b[$i])
    $styleText.Copy()
    $c.PasteSpecial(-4122)
    $c.Value = $values2b[$i]
}

Write-Host "B2C hotel booking sheet + flight employee columns added"
